$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_QueryLog_External")
Write-Host "Sheet found"
Write-Host $ws.Range("A2").Value
